$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New block header (H1:J1): merge first, then mirror the BMW block (E1:G1) format ---
$ws.Range("H1:J1").Merge()
$ws.Range("H1").Value2 = "Bosch 0280155746 200cc"
$ws.Range("E1:G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)   # xlPasteFormats

# --- New block sub-header (H2:J2): ON / OFF / LAG, mirrors E2:G2 ---
$ws.Range("H2").Value2 = "ON"
$ws.Range("I2").Value2 = "OFF"
$ws.Range("J2").Value2 = "LAG"
$ws.Range("E2:G2").Copy()
$ws.Range("H2:J2").PasteSpecial(-4122)

# --- Data rows 3-8: values + H-I formula, formatting mirrored row-by-row from E:G ---
$ws.Range("H3").Value2 = 2.71
$ws.Range("I3").Value2 = 0.85
$ws.Range("J3").Formula = "=H3-I3"
$ws.Range("E3:G3").Copy()
$ws.Range("H3:J3").PasteSpecial(-4122)

$ws.Range("H4").Value2 = 1.93
$ws.Range("I4").Value2 = 0.85
$ws.Range("J4").Formula = "=H4-I4"
$ws.Range("E4:G4").Copy()
$ws.Range("H4:J4").PasteSpecial(-4122)

$ws.Range("H5").Value2 = 1.62
$ws.Range("I5").Value2 = 0.85
$ws.Range("J5").Formula = "=H5-I5"
$ws.Range("E5:G5").Copy()
$ws.Range("H5:J5").PasteSpecial(-4122)

$ws.Range("H6").Value2 = 1.32
$ws.Range("I6").Value2 = 0.86
$ws.Range("J6").Formula = "=H6-I6"
$ws.Range("E6:G6").Copy()
$ws.Range("H6:J6").PasteSpecial(-4122)

$ws.Range("H7").Value2 = 1.1599999999999999
$ws.Range("I7").Value2 = 0.87
$ws.Range("J7").Formula = "=H7-I7"
$ws.Range("E7:G7").Copy()
$ws.Range("H7:J7").PasteSpecial(-4122)

$ws.Range("H8").Value2 = 0.99
$ws.Range("I8").Value2 = 0.88
$ws.Range("J8").Formula = "=H8-I8"
$ws.Range("E8:G8").Copy()
$ws.Range("H8:J8").PasteSpecial(-4122)

# --- View state: zoom + selection, matching the author's final view ---
$excel.ActiveWindow.Zoom = 160
$ws.Range("P27").Select() | Out-Null

